# Finance - Core Banking System Modernization KPI Dashboard
# Rebuild AI/ML sample content as Finance / Core Banking System Modernization content.

$wb = $excel.ActiveWorkbook

$wsInstructions = $wb.Worksheets.Item("Instructions & User Guide")
$wsDashboard    = $wb.Worksheets.Item("KPI Dashboard")

# ---------------------------------------------------------------------------
# Sheet 1: "Instructions & User Guide" - title banner
# ---------------------------------------------------------------------------
$wsInstructions.Range("A1").Value = "Finance - Core Banking System Modernization KPI Dashboard - User Guide & Instructions"

# Re-materialize the blank spacer rows between sections (rows are present but
# empty in the source template; touching a no-op row property keeps them
# explicit in the saved sheet without altering their formatting).
foreach ($r in @(9, 16, 25, 33, 40, 47, 56, 57)) {
    $wsInstructions.Rows.Item($r).OutlineLevel = 0
}

# ---------------------------------------------------------------------------
# Sheet 2: "KPI Dashboard" - header banner
# ---------------------------------------------------------------------------
$wsDashboard.Range("A1").Value = "FINANCE - CORE BANKING SYSTEM MODERNIZATION - KPI DASHBOARD"
$wsDashboard.Range("A2").Value = "Project: Core Banking System Modernization"

# Re-materialize the blank spacer row above the KPI table header.
$wsDashboard.Rows.Item(5).OutlineLevel = 0

# ---------------------------------------------------------------------------
# KPI table (rows 8-22): name / category / target / current / status, plus the
# Variance_% column which switches from a static percent-string to a live
# formula, and the Notes column which is reworded for the new project.
# ---------------------------------------------------------------------------

$noteText = "Critical KPI for Finance - Core Banking System Modernization success"

$kpis = @(
    @{ Row = 8;  B = "Transaction Processing Accuracy"; C = "Performance"; D = 99.95;    E = 99.87;  G = "At Risk" }
    @{ Row = 9;  B = "System Uptime %";                 C = "Performance"; D = 99.9;      E = 99.92;  G = "On Track" }
    @{ Row = 10; B = "Data Migration Completion";       C = "Performance"; D = 95;        E = 88;     G = "At Risk" }
    @{ Row = 11; B = "User Adoption Rate";              C = "Performance"; D = 85;        E = 78;     G = "At Risk" }
    @{ Row = 12; B = "Regulatory Compliance Score";     C = "Quality";     D = 100;       E = 98;     G = "On Track" }
    @{ Row = 13; B = "Security Incident Count";         C = "Quality";     D = 0;         E = 2;      G = "At Risk" }
    @{ Row = 14; B = "Defect Density";                  C = "Quality";     D = 0.5;       E = 0.8;    G = "At Risk" }
    @{ Row = 15; B = "Training Completion Rate";        C = "Quality";     D = 95;        E = 92;     G = "On Track" }
    @{ Row = 16; B = "Budget Variance %";                C = "Financial";   D = 0;         E = 3.5;    G = "At Risk" }
    @{ Row = 17; B = "Schedule Variance %";              C = "Financial";   D = 0;         E = 2.8;    G = "At Risk" }
    @{ Row = 18; B = "Cost Savings Achieved";            C = "Financial";   D = 2500000;   E = 2100000; G = "At Risk" }
    @{ Row = 19; B = "ROI Achievement %";                C = "Financial";   D = 125;       E = 118;    G = "At Risk" }
    @{ Row = 20; B = "Customer Satisfaction Score";      C = "Financial";   D = 4.5;       E = 4.2;    G = "On Track" }
    @{ Row = 21; B = "Integration Success Rate";         C = "Financial";   D = 98;        E = 95;     G = "On Track" }
    @{ Row = 22; B = "Performance Benchmark Achievement"; C = "Financial";  D = 100;       E = 92;     G = "At Risk" }
)

foreach ($kpi in $kpis) {
    $r = $kpi.Row

    $wsDashboard.Range("B$r").Value = $kpi.B
    $wsDashboard.Range("C$r").Value = $kpi.C
    $wsDashboard.Range("D$r").Value = $kpi.D
    $wsDashboard.Range("E$r").Value = $kpi.E
    $wsDashboard.Range("F$r").Formula = "=((E$r-D$r)/D$r)*100"
    $wsDashboard.Range("G$r").Value = $kpi.G
    $wsDashboard.Range("K$r").Value = $noteText
}

Write-Host "Finance KPI Dashboard content applied."
